# Applies the "Skrevet litt i Prosjektbeskrivelsen" edit to the project
# description document.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Heading: "tittel" -> "Musikkavisen tungrocken"
# ------------------------------------------------------------------
$d.Content.Find.Execute("tittel", $true, $false, $false, $false, $false, $true, 1, $false, "Musikkavisen tungrocken", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Body paragraph: "Tekst" -> the first sentence, followed by the
#    rest of the paragraph appended as additional text.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Tekst", $true, $false, $false, $false, $false, $true, 1, $false, "Gjennom faget ", 2) | Out-Null

$introPara = $d.Paragraphs(8)
$pos = $introPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Datamodellering og databaseapplikasjoner")

$introPara = $d.Paragraphs(8)
$pos = $introPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter(" i et tidligere semester, lagde vi som gruppeprosjekt et CMS (Content Management System) for en fiktiv musikkavis med navnet Tungrocken. I mangel på bedre og mer kreative ideer, så har vi tenkt å bygge denne ")

$introPara = $d.Paragraphs(8)
$pos = $introPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter("web-baserte nyhetstjenesten ")

$introPara = $d.Paragraphs(8)
$pos = $introPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter("om til en app-løsning, med tilhørende server for datahåndtering. ")

# ------------------------------------------------------------------
# 3. Insert four brand-new paragraphs after the intro paragraph.
# ------------------------------------------------------------------
$introPara = $d.Paragraphs(8)
$introPara.Range.InsertParagraphAfter() | Out-Null

$restPara = $d.Paragraphs(9)
$restPara.Range.Text = "Vi vil forsøke å sette opp en REST-basert serverløsning ved hjelp av JAX-RS"
$restPara = $d.Paragraphs(9)
$pos = $restPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter(".")

$restPara = $d.Paragraphs(9)
$restPara.Range.InsertParagraphAfter() | Out-Null

$appPara = $d.Paragraphs(10)
$appPara.Range.Text = "App"
$appPara = $d.Paragraphs(10)
$pos = $appPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter("´")
$appPara = $d.Paragraphs(10)
$pos = $appPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter("en vil bli satt opp som en abonnement-tjeneste, der du må være en registrert bruker for å få tilgang til innholdet. ")
$appPara = $d.Paragraphs(10)
$pos = $appPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Om")
$appPara = $d.Paragraphs(10)
$pos = $appPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter(" en bruker ikke er logget inn, vil et vindu for pålogging dukke opp. Der vil det også være mulig å klikke seg videre for registrering av ny bruker.")

$appPara = $d.Paragraphs(10)
$appPara.Range.InsertParagraphAfter() | Out-Null

$loginPara = $d.Paragraphs(11)
$loginPara.Range.Text = "Når brukeren er innlogget, og har fått nødvendig tilgang"
$loginPara = $d.Paragraphs(11)
$pos = $loginPara.Range.End - 1
$d.Range($pos, $pos).InsertAfter(", vil hovedsiden dukke opp. Hovedsiden vil vise en liten smakebit av nyhetsartiklene, som det er mulig å klikke seg videre inn på for å få opp hele (altså som en vanlig nyhetsapp…)")

$loginPara = $d.Paragraphs(11)
$loginPara.Range.InsertParagraphAfter() | Out-Null

$layoutPara = $d.Paragraphs(12)
$layoutPara.Range.Text = "Oppsettet i hovedsiden vil være slik at den nyeste artikkelen vil vise i større format enn de resterende, som vil bli listet opp kontinuerlig. "

# ------------------------------------------------------------------
# 4. The final, originally-empty paragraph holds the "_GoBack"
#    bookmark. Type the admin-panel text into it *before* the
#    bookmark's position so the bookmark ends up after the text
#    (matching how Word leaves _GoBack at the last edit point).
# ------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs(13)
$startPos = $bookmarkPara.Range.Start
$d.Range($startPos, $startPos).InsertBefore("I utgangspunktet vil vi også lage til et kontrollpanel for administrator, som kan ")

$bookmarkPara = $d.Paragraphs(13)
$pos = $bookmarkPara.Range.End - 1
$d.Range($pos, $pos).InsertBefore("legge til rette for redigering av brukere og innhold (legge til artikler, endre artikler, deaktivere artikler, endre brukere, deaktivere brukere).")

# ------------------------------------------------------------------
# 5. Add a new, empty trailing paragraph after the bookmark paragraph.
# ------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs(13)
$bookmarkPara.Range.InsertParagraphAfter() | Out-Null
